$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text cell A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.12 = 7725.63 pesos`n✅ 7725.63 pesos = 2.11 = 938.04 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the N10, O10 and O12 values ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 470.777
$tasas.Range("O10").Value = 3637.05
$tasas.Range("O12").Value = 445
